## Update gh-pages to output generated at 456a3b4
## Applies numeric "want-to-go" count bumps across the four sheets, plus the
## weekly reshuffle of rows 22-26 on the "展览" sheet (一场 cancelled show
## drops out, the remaining shows shift up one row, and a new show is
## appended at the bottom of the block).

$wb = $excel.ActiveWorkbook

function Set-Num($ws, $row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.Value = $val
}

function Set-Text($ws, $row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.Value = $val
}

# Some "date-looking" text cells (B column, e.g. "2024-10-04") would be
# auto-converted to real Excel dates on assignment; force them back to
# plain text by flipping the format to Text for the write, then clearing
# the format again so no stray style sticks around on the cell.
function Set-TextForceString($ws, $row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

# ---------------------------------------------------------------------
# Sheet 1: 展览 (exhibitions)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

Set-Num $ws1 4  6 496
Set-Num $ws1 5  6 8944
Set-Num $ws1 6  6 8944
Set-Num $ws1 8  6 11495
Set-Num $ws1 14 6 129

# Row 22 <- old row 23 content (北京·美漫超级英雄ONLY（取消） drops out)
Set-TextForceString $ws1 22 2 "2024-10-04"
Set-Text $ws1 22 3 "北京·第五人格only同人展"
Set-Text $ws1 22 4 "北花园路1号超级蜂巢C座 超级蜂巢国际会议中心"
Set-Text $ws1 22 5 "2024.10.04 10:00-10.04 17:00"
Set-Num  $ws1 22 6 1956
Set-Num  $ws1 22 7 68
Set-Text $ws1 22 8 "https://show.bilibili.com/platform/detail.html?id=89309"
Set-Text $ws1 22 9 "//i0.hdslb.com/bfs/openplatform/202407/4XsICpa71721046044404.jpeg"

# Row 23 <- old row 24 content
Set-Text $ws1 23 3 "北京·首届SH动漫游戏展"
Set-Text $ws1 23 4 "安定路5号院(安贞门地铁站A西北口步行420米) 北京北投购物公园"
Set-Text $ws1 23 5 "2024.10.04 10:00-10.04 18:00"
Set-Num  $ws1 23 6 751
Set-Num  $ws1 23 7 55
Set-Text $ws1 23 8 "https://show.bilibili.com/platform/detail.html?id=91635"
Set-Text $ws1 23 9 "//i2.hdslb.com/bfs/openplatform/202409/SKe1HMLH1725179765551.jpeg"

# Row 24 <- old row 25 content
Set-Text $ws1 24 3 "帝都·重返未来1999同人ONLY金秋深眠"
Set-Text $ws1 24 4 "华佗路与新源大街交汇处西100米 凯德MALL·大兴"
Set-Text $ws1 24 5 "2024.10.04 10:00-10.05 17:00"
Set-Num  $ws1 24 6 689
Set-Num  $ws1 24 7 68
Set-Text $ws1 24 8 "https://show.bilibili.com/platform/detail.html?id=92315"
Set-Text $ws1 24 9 "//i1.hdslb.com/bfs/openplatform/202409/YHMYHehz1726129707544.jpeg"

# Row 25 <- old row 26 content
Set-TextForceString $ws1 25 2 "2024-10-05"
Set-Text $ws1 25 3 "北京·咒术回战同人Only2.0"
Set-Text $ws1 25 4 "安定路5号院(安贞门地铁站A西北口步行420米) 北京北投购物公园"
Set-Text $ws1 25 5 "2024.10.05 09:30-10.05 17:00"
Set-Num  $ws1 25 6 372
Set-Num  $ws1 25 7 65
Set-Text $ws1 25 8 "https://show.bilibili.com/platform/detail.html?id=91628"
Set-Text $ws1 25 9 "//i0.hdslb.com/bfs/openplatform/202408/IsJo7aU61724405528082.jpeg"

# Row 26 <- brand-new show (北京·马娘同人ONLY2 in row 27 is untouched)
Set-Text $ws1 26 3 "北京·悠唐购物中心国庆动漫嘉年华（免费）"
Set-Text $ws1 26 4 "悠唐购物中心 悠唐购物中心"
Set-Text $ws1 26 5 "2024.10.05 13:00-10.05 18:00"
Set-Num  $ws1 26 6 0
Set-Num  $ws1 26 7 30
Set-Text $ws1 26 8 "https://show.bilibili.com/platform/detail.html?id=93026"
Set-Text $ws1 26 9 "//i1.hdslb.com/bfs/openplatform/202409/vVftaNeG1727434341484.jpeg"

Set-Num $ws1 29 6 623
Set-Num $ws1 30 6 1400
Set-Num $ws1 34 6 48
Set-Num $ws1 38 6 319
Set-Num $ws1 39 6 90
Set-Num $ws1 40 6 364
Set-Num $ws1 42 6 150
Set-Num $ws1 45 6 136
Set-Num $ws1 46 6 836

# ---------------------------------------------------------------------
# Sheet 2: 演出 (performances)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
Set-Num $ws2 14 6 39
Set-Num $ws2 24 6 89

# ---------------------------------------------------------------------
# Sheet 3: 本地生活 (local life)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
Set-Num $ws3 3 6 2877

# ---------------------------------------------------------------------
# Sheet 4: 全部类型 (all categories) -- same underlying events as sheet 1
# but at different row offsets; only the "want-to-go" counters move here,
# the cancelled/added rows from sheet 1 are not represented on this sheet.
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
Set-Num $ws4 6  6 496
Set-Num $ws4 7  6 8944
Set-Num $ws4 8  6 8944
Set-Num $ws4 10 6 11495
Set-Num $ws4 15 6 129
Set-Num $ws4 19 6 1956
Set-Num $ws4 20 6 751
Set-Num $ws4 21 6 689
Set-Num $ws4 22 6 372
Set-Num $ws4 25 6 623
Set-Num $ws4 28 6 1400
Set-Num $ws4 33 6 48
Set-Num $ws4 36 6 319
Set-Num $ws4 38 6 364
Set-Num $ws4 40 6 150
Set-Num $ws4 43 6 136

Write-Output "done"
